$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores figures as literal text in the source data
# (values like "46.756.12" or "0.0800" are not valid single numbers /
# would lose trailing zeros if Excel auto-converted them). Force every
# Price cell we touch to Text format first so the assigned string is
# kept verbatim instead of being reinterpreted as a number.
foreach ($addr in @("D2","D3","D5","D6","D7","D9","D10","D11","D12","D14","D15","D16","D18","D19","D20","D21","D22","D23","D27","D29","D30","D31","D32","D39","D40","D41","D43","D44","D45","D46","D47","D48","D49","D50","D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "46.756.12"
$ws.Range("E2").Value = "  +6.47%  "

# Row 3
$ws.Range("D3").Value = "2.314.86"
$ws.Range("E3").Value = "  +5.36%  "

# Row 4
$ws.Range("E4").Value = "  -0.56%  "

# Row 5
$ws.Range("D5").Value = "300.28"
$ws.Range("E5").Value = "  +2.17%  "

# Row 6
$ws.Range("D6").Value = "98.96"
$ws.Range("E6").Value = "  +10.81%  "

# Row 7
$ws.Range("D7").Value = "0.574"
$ws.Range("E7").Value = "  +0.44%  "

# Row 8
$ws.Range("E8").Value = "  -0.47%  "

# Row 9
$ws.Range("D9").Value = "0.531"
$ws.Range("E9").Value = "  +10.64%  "

# Row 10
$ws.Range("D10").Value = "35.63"
$ws.Range("E10").Value = "  +9.85%  "

# Row 11
$ws.Range("D11").Value = "0.0800"
$ws.Range("E11").Value = "  +3.50%  "

# Row 12
$ws.Range("D12").Value = "7.37"
$ws.Range("E12").Value = "  +9.44%  "

# Row 13
$ws.Range("E13").Value = "  +1.67%  "

# Row 14
$ws.Range("D14").Value = "2.666.69"
$ws.Range("E14").Value = "  +5.33%  "

# Row 15
$ws.Range("D15").Value = "2.312.37"
$ws.Range("E15").Value = "  +1.88%  "

# Row 16
$ws.Range("D16").Value = "13.98"
$ws.Range("E16").Value = "  +7.17%  "

# Row 17
$ws.Range("E17").Value = "  +7.47%  "

# Row 18
$ws.Range("D18").Value = "46.698.22"
$ws.Range("E18").Value = "  +6.49%  "

# Row 19
$ws.Range("D19").Value = "13.36"
$ws.Range("E19").Value = "  +23.95%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0942"
$ws.Range("E20").Value = "  +6.85%  "

# Row 21
$ws.Range("D21").Value = "6.16"
$ws.Range("E21").Value = "  +5.58%  "

# Row 22
$ws.Range("D22").Value = "66.93"
$ws.Range("E22").Value = "  +6.22%  "

# Row 23
$ws.Range("D23").Value = "249.48"
$ws.Range("E23").Value = "  +9.16%  "

# Row 24
$ws.Range("E24").Value = "  +6.38%  "

# Row 25
$ws.Range("E25").Value = "  +9.70%  "

# Row 26
$ws.Range("E26").Value = "  +0.00%  "

# Row 27
$ws.Range("D27").Value = "43.11"
$ws.Range("E27").Value = "  +21.70%  "

# Row 28
$ws.Range("E28").Value = "  +1.73%  "

# Row 29
$ws.Range("D29").Value = "9.88"
$ws.Range("E29").Value = "  +7.68%  "

# Row 30
$ws.Range("D30").Value = "20.14"
$ws.Range("E30").Value = "  +6.39%  "

# Row 31
$ws.Range("D31").Value = "5.81"
$ws.Range("E31").Value = "  +9.61%  "

# Row 32
$ws.Range("D32").Value = "147.79"
$ws.Range("E32").Value = "  -0.59%  "

# Row 33
$ws.Range("E33").Value = "  +9.75%  "

# Row 34
$ws.Range("E34").Value = "  +6.05%  "

# Row 35
$ws.Range("E35").Value = "  +9.75%  "

# Row 36
$ws.Range("E36").Value = "  +10.77%  "

# Row 37
$ws.Range("E37").Value = "  +2.87%  "

# Row 38
$ws.Range("E38").Value = "  +10.56%  "

# Row 39
$ws.Range("D39").Value = "15.63"
$ws.Range("E39").Value = "  +17.26%  "

# Row 40
$ws.Range("D40").Value = "4.01"
$ws.Range("E40").Value = "  +15.34%  "

# Row 41
$ws.Range("D41").Value = "3.47"
$ws.Range("E41").Value = "  +13.08%  "

# Row 42
$ws.Range("E42").Value = "  +10.14%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "2.01"
$ws.Range("E43").Value = "  +22.49%  "

# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.73%  "

# Row 45
$ws.Range("D45").Value = "1.844.22"
$ws.Range("E45").Value = "  +6.11%  "

# Row 46
$ws.Range("D46").Value = "90.87"
$ws.Range("E46").Value = "  +23.53%  "

# Row 47
$ws.Range("D47").Value = "0.201"
$ws.Range("E47").Value = "  +17.38%  "

# Row 48
$ws.Range("D48").Value = "75.51"
$ws.Range("E48").Value = "  +11.27%  "

# Row 49
$ws.Range("D49").Value = "4.96"
$ws.Range("E49").Value = "  +12.79%  "

# Row 50
$ws.Range("D50").Value = "97.50"
$ws.Range("E50").Value = "  +6.85%  "

# Row 51
$ws.Range("D51").Value = "54.31"
$ws.Range("E51").Value = "  +10.66%  "
